$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old, unused leading column (A) which held duplicate GENE values
# with header styling but no header text. This shifts B:F left to A:E.
$ws.Columns("A").Delete()

# Fix header typo: MODEL_CONDITION -> MODELCONDITION (now in column D after shift)
$ws.Cells.Item(1, 4).Value = "MODELCONDITION"
